# Update MAIN DASHBOARD station load values (columns B,C,D,F,H,I) for hours 1-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Hour 1)
$ws.Range("B2").Value = 22266
$ws.Range("C2").Value = 5999
$ws.Range("D2").Value = 9984
$ws.Range("F2").Value = 6030
$ws.Range("H2").Value = 17786
$ws.Range("I2").Value = 16906

# Row 3 (Hour 2)
$ws.Range("B3").Value = 21326
$ws.Range("C3").Value = 5717
$ws.Range("D3").Value = 9532
$ws.Range("F3").Value = 5923
$ws.Range("H3").Value = 17126
$ws.Range("I3").Value = 16079

# Row 4 (Hour 3)
$ws.Range("B4").Value = 20449
$ws.Range("C4").Value = 5454
$ws.Range("D4").Value = 9230
$ws.Range("F4").Value = 5604
$ws.Range("H4").Value = 16410
$ws.Range("I4").Value = 15490

# Row 5 (Hour 4)
$ws.Range("B5").Value = 19762
$ws.Range("C5").Value = 5204
$ws.Range("D5").Value = 8913
$ws.Range("F5").Value = 5557
$ws.Range("H5").Value = 15487
$ws.Range("I5").Value = 15148

# Row 6 (Hour 5)
$ws.Range("B6").Value = 19910
$ws.Range("C6").Value = 5266
$ws.Range("D6").Value = 8941
$ws.Range("F6").Value = 5670
$ws.Range("H6").Value = 15633
$ws.Range("I6").Value = 15424

# Row 7 (Hour 6)
$ws.Range("B7").Value = 19707
$ws.Range("C7").Value = 5468
$ws.Range("D7").Value = 9267
$ws.Range("F7").Value = 7157
$ws.Range("H7").Value = 16091
$ws.Range("I7").Value = 15868

# Row 8 (Hour 7)
$ws.Range("B8").Value = 20268
$ws.Range("C8").Value = 4604
$ws.Range("D8").Value = 9666
$ws.Range("F8").Value = 7401
$ws.Range("H8").Value = 15557
$ws.Range("I8").Value = 16013

# Row 9 (Hour 8)
$ws.Range("B9").Value = 23185
$ws.Range("C9").Value = 4901
$ws.Range("D9").Value = 11223
$ws.Range("F9").Value = 10623
$ws.Range("H9").Value = 18375
$ws.Range("I9").Value = 18064

# Row 10 (Hour 9)
$ws.Range("B10").Value = 28823
$ws.Range("C10").Value = 5818
$ws.Range("D10").Value = 13930
$ws.Range("F10").Value = 14092
$ws.Range("H10").Value = 20118
$ws.Range("I10").Value = 21041

# Row 11 (Hour 10)
$ws.Range("B11").Value = 32025
$ws.Range("C11").Value = 6084
$ws.Range("D11").Value = 18165
$ws.Range("F11").Value = 15134
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 22980

# Row 12 (Hour 11)
$ws.Range("B12").Value = 32138
$ws.Range("C12").Value = 6317
$ws.Range("D12").Value = 19577
$ws.Range("F12").Value = 15265
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 23801

# Row 13 (Hour 12)
$ws.Range("B13").Value = 31865
$ws.Range("C13").Value = 6389
$ws.Range("D13").Value = 19812
$ws.Range("F13").Value = 15463
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 23781

# Row 14 (Hour 13)
$ws.Range("B14").Value = 32062
$ws.Range("C14").Value = 6461
$ws.Range("D14").Value = 19870
$ws.Range("F14").Value = 15114
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 23508

# Row 15 (Hour 14)
$ws.Range("B15").Value = 33756
$ws.Range("C15").Value = 6996
$ws.Range("D15").Value = 20143
$ws.Range("F15").Value = 16097
$ws.Range("H15").Value = 22648
$ws.Range("I15").Value = 24132

# Row 16 (Hour 15)
$ws.Range("B16").Value = 33340
$ws.Range("C16").Value = 6897
$ws.Range("D16").Value = 19871
$ws.Range("F16").Value = 15846
$ws.Range("H16").Value = 22412
$ws.Range("I16").Value = 24607

# Row 17 (Hour 16)
$ws.Range("B17").Value = 33211
$ws.Range("C17").Value = 6639
$ws.Range("D17").Value = 20137
$ws.Range("H17").Value = 22106
$ws.Range("I17").Value = 0

# Row 18 (Hour 17)
$ws.Range("B18").Value = 30925
$ws.Range("C18").Value = 6519
$ws.Range("D18").Value = 19370
